# "version final sin errores"
#
# Changes applied to the "Metadata" sheet:
#   1. Bump the Version value (B3) from "0.4.0" to "0.7.0".
#   2. Remove the "Jurisdiction" / "Chile" row entirely (row 11), which
#      shifts every row below it up by one (old A12:B15 -> new A11:B14)
#      and shrinks the sheet's used range from A1:B15 to A1:B14.
#
# Sheet2 ("Include from ...") is left untouched content-wise; its shared
# string indices simply follow the shrink of the shared string table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# 1. Update Version: 0.4.0 -> 0.7.0
$ws.Range("B3").Value = "0.7.0"

# 2. Delete the whole "Jurisdiction" / "Chile" row (row 11)
$ws.Range("A11").EntireRow.Delete()
